$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new trade row (row 7)
$ws.Range("A7").Value = 42650.366909722223
$ws.Range("A7").NumberFormat = "m/d/yy h:mm"
$ws.Range("B7").Value = $false
$ws.Range("C7").Value = 9814.67
$ws.Range("D7").Value = 9852.6
$ws.Range("E7").Value = 104.839996
$ws.Range("F7").Value = 104.029999
$ws.Range("G7").Value = $false
$ws.Range("H7").Value = -0.77
$ws.Range("I7").Value = $false

# Column E's longest value is now as wide as column F's ("104.839996" matches
# "104.029999"/"105.290001"), so Excel's best-fit recalculation widens E to
# match F - the two adjacent same-width columns collapse into one <col> span.
$ws.Range("E1:F1").ColumnWidth = 10
